$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text looks numeric: force Text format so Excel
# does not auto-convert the literal string into a number (losing
# formatting like trailing/leading zeros).
$textCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D15", "D16", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D39", "D43", "D45", "D46", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.448.58'
$ws.Range("D3").Value = '2.377.15'
$ws.Range("E3").Value = '  +4.89%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '235.39'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").Value = '0.652'
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = '72.32'
$ws.Range("E7").Value = '  +13.71%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +5.19%  '
$ws.Range("D10").Value = '0.0980'
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("D11").Value = '56.85'
$ws.Range("E11").Value = '  -2.23%  '
$ws.Range("D12").Value = '27.25'
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("D13").Value = '2.734.46'
$ws.Range("E13").Value = '  +5.04%  '
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").Value = '16.05'
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").Value = '6.31'
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("D17").Value = '0.858'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").Value = '2.377.19'
$ws.Range("E18").Value = '  +4.79%  '
$ws.Range("D19").Value = '43.446.07'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").Value = '0.0000100'
$ws.Range("E20").Value = '  +1.68%  '
$ws.Range("D21").Value = '6.36'
$ws.Range("E21").Value = '  +3.06%  '
$ws.Range("D22").Value = '74.65'
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").Value = '251.12'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '3.63'
$ws.Range("E26").Value = '  +8.15%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '10.06'
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").Value = '22.63'
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("D30").Value = '174.53'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").Value = '1.52'
$ws.Range("E31").Value = '  +5.85%  '
$ws.Range("E32").Value = '  -5.08%  '
$ws.Range("D33").Value = '0.127'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = '5.02'
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = '0.0693'
$ws.Range("E35").Value = '  +0.92%  '
$ws.Range("D36").Value = '5.10'
$ws.Range("E36").Value = '  +2.66%  '
$ws.Range("E37").Value = '  +7.16%  '
$ws.Range("E38").Value = '  +2.41%  '
$ws.Range("D39").Value = '3.68'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("E42").Value = '  +2.10%  '
$ws.Range("D43").Value = '18.60'
$ws.Range("E43").Value = '  +7.71%  '
$ws.Range("E44").Value = '  +8.78%  '
$ws.Range("D45").Value = '100.65'
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("D46").Value = '4.51'
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("E47").Value = '  +2.14%  '
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.604.64'
$ws.Range("E50").Value = '  +5.16%  '
$ws.Range("B51").Value = 'TerraClassic'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D51").Value = '0.000205'
$ws.Range("E51").Value = '  -7.31%  '

# Restore default style on those cells (keeps cell-level style index
# identical to the original file; only the underlying text changed).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
